$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Slurry")

# Update digestate (DM) column C for rows 2-5: 5.1% -> 5.9%
$ws.Range("C2").Value = 5.9
$ws.Range("C3").Value = 5.9
$ws.Range("C4").Value = 5.9
$ws.Range("C5").Value = 5.9

# Remove the now-duplicate rows 6-9 (the 6.9% block), leaving only the 5.9% rows
$ws.Rows("6:9").Delete()

$ws.Range("F14").Select()
